$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 141 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
